$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old helper/unused rows (original sheet had data through row 9,
# the updated sheet only has data through row 5).
$ws.Rows("6:9").Delete()

# Column C previously held a formula ("=A2" etc, shared down the column).
# It is now a plain picked-winner value, so clear any leftover formulas
# before writing the new static values.
$ws.Range("C2:C5").ClearContents()

# Row 2: Auburn vs Michigan St. -> Auburn wins
$ws.Range("A2").Value = "Auburn"
$ws.Range("B2").Value = "Michigan St."
$ws.Range("C2").Value = "Auburn"
$ws.Range("D2").Value = 5.8595544965272603

# Row 3: Florida vs Texas Tech -> Texas Tech wins
$ws.Range("A3").Value = "Florida"
$ws.Range("B3").Value = "Texas Tech"
$ws.Range("C3").Value = "Texas Tech"
$ws.Range("D3").Value = 0.40823675122630698

# Row 4: Duke vs Alabama -> Duke wins
$ws.Range("A4").Value = "Duke"
$ws.Range("B4").Value = "Alabama"
$ws.Range("C4").Value = "Duke"
$ws.Range("D4").Value = 1.956564658704868

# Row 5: Houston vs Tennessee -> Houston wins
$ws.Range("A5").Value = "Houston"
$ws.Range("B5").Value = "Tennessee"
$ws.Range("C5").Value = "Houston"
$ws.Range("D5").Value = 0.85560703208208611

$ws.Range("D2").Select()

$wb.Save()
